$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = -10.534
    32  = -13.239
    36  = -12.732
    38  = -12.544
    46  = -13.935
    54  = -12.737
    55  = -13.65
    67  = -11.518
    69  = -10.744
    72  = -11.753
    91  = -12.228
    99  = -11.412
    104 = -12.729
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
